$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "Design"
$design = $wb.Worksheets.Item(1)
$design.Name = "Design"

# Add a new "Questions" sheet after Design
$questions = $wb.Worksheets.Add($null, $design)
$questions.Name = "Questions"

$questions.Range("A5").Value = "* Is there a PostgreSQL command that returns the primary key and foreign keys for a table?"
$questions.Range("A6").Value = "* How do we handle foreign key columns that have different meaning in a table?"
$questions.Range("B7").Value = "There can be multiple versions of the same foreign key with different meanings."
$questions.Range("A8").Value = "* If we split all data into multiple tables, how can we make a table that has the values of the keys, not the integer keys themselves?"
$questions.Range("B9").Value = "Is it a matter of doing a series of joins? That seems tedious."

$questions.Activate()
